# Adding the changes we made on may 9th
# Insert 11 new sensor-reading rows right after the header row (pushing the
# existing data down), and drop the oldest reading that now falls off the
# end of the tracked window (old row 21 -> row 32 after the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer samples to insert at the top of the data (rows 2-12)
$newData = @(
    @(-4.1851313829422, 4.315697193145752, -1.521726083755494),
    @(-4.044641017913818, 4.278119564056396, -1.635269999504089),
    @(-4.037833023071289, 4.246780538558959, -1.5330271422863),
    @(-4.188077163696289, 4.206779670715332, -1.555011582374573),
    @(-4.266698265075684, 4.172728943824768, -1.657637816667557),
    @(-4.307723093032837, 4.115167903900146, -1.525768172740936),
    @(-4.275023818016052, 4.170393228530884, -1.397946566343308),
    @(-4.101493644714355, 4.221428394317627, -1.457295727729798),
    @(-4.052593421936035, 4.150938177108764, -1.483953177928925),
    @(-4.146867895126342, 4.078793859481812, -1.569497275352478),
    @(-4.033952045440674, 4.105497026443481, -1.389435482025145)
)

$insertCount = $newData.Count

# Shift the existing data rows (2..21) down by inserting new rows above row 2
$insertRange = $ws.Range("A2:C$($insertCount + 1)")
$insertRange.Insert()
# Inserting copies the header row's formatting down onto the new rows;
# strip it so the new data cells stay unstyled like the rest of the data.
$insertRange.ClearFormats()

# Write the new rows into the freshly inserted space
for ($i = 0; $i -lt $insertCount; $i++) {
    $rowNum = 2 + $i
    $vals = $newData[$i]
    $ws.Cells.Item($rowNum, 1).Value = $vals[0]
    $ws.Cells.Item($rowNum, 2).Value = $vals[1]
    $ws.Cells.Item($rowNum, 3).Value = $vals[2]
}

# The window keeps only 30 data rows; drop the oldest sample that now
# falls past row 31 (original row 21, now shifted to row 32).
$ws.Rows.Item(32).Delete()
